$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @(new D value or $null, new E value or $null, D-is-numeric-looking)
$updates = @{
    2 = @("29.160.51", "  -0.16%  ", $false)
    3 = @("1.837.80", "  -0.57%  ", $false)
    4 = @($null, "  -0.01%  ", $false)
    5 = @("240.36", "  -2.41%  ", $true)
    6 = @("0.6858", "  -1.86%  ", $true)
    7 = @("0.9993", "  -0.01%  ", $true)
    8 = @($null, "  -1.68%  ", $false)
    9 = @("0.07463", "  -3.59%  ", $true)
    10 = @("23.16", "  -1.75%  ", $true)
    11 = @("0.07665", "  -2.00%  ", $true)
    12 = @("1.838.46", "  -0.39%  ", $false)
    13 = @("5.059", "  -1.42%  ", $true)
    14 = @("0.6820", "  -0.72%  ", $true)
    15 = @("87.63", "  -6.35%  ", $true)
    16 = @("6.165", "  -7.26%  ", $true)
    17 = @("29.138.16", "  -0.16%  ", $false)
    18 = @("0.000008169", "  -1.92%  ", $true)
    19 = @("2.081.07", "  -0.08%  ", $false)
    20 = @("227.84", "  -5.66%  ", $true)
    21 = @("12.53", "  -1.94%  ", $true)
    22 = @($null, "  +0.00%  ", $false)
    23 = @("7.401", $null, $true)
    24 = @("0.9993", "  -0.02%  ", $true)
    25 = @("0.1456", "  -4.08%  ", $true)
    26 = @("160.02", "  +0.61%  ", $true)
    27 = @("8.763", $null, $true)
    28 = @($null, "  -1.00%  ", $false)
    29 = @("1.513", "  -2.00%  ", $true)
    30 = @("4.271", "  +0.77%  ", $true)
    31 = @("4.151", "  -1.04%  ", $true)
    32 = @($null, "  +0.33%  ", $false)
    33 = @("0.05176", "  +1.03%  ", $true)
    34 = @("0.7663", "  -3.24%  ", $true)
    35 = @("1.844", "  -1.30%  ", $true)
    36 = @("1.135", "  -1.40%  ", $true)
    37 = @($null, "  -0.59%  ", $false)
    38 = @("1.313.70", "  +0.20%  ", $false)
    39 = @("0.01834", "  -1.96%  ", $true)
    40 = @("2.721", "  +0.36%  ", $true)
    41 = @("0.9346", "  -1.27%  ", $true)
    42 = @("5.793", "  -4.72%  ", $true)
    43 = @("104.70", "  -2.77%  ", $true)
    44 = @("0.9986", "  -0.07%  ", $true)
    45 = @($null, "  +0.22%  ", $false)
    46 = @("65.14", "  +1.39%  ", $true)
    47 = @("1.981.91", $null, $false)
    48 = @("0.5202", "  +0.45%  ", $true)
    49 = @("9.531", "  -2.09%  ", $true)
    50 = @($null, "  +0.20%  ", $false)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    $dIsNumericLooking = $vals[2]
    if ($dVal -ne $null) {
        $dCell = $ws.Cells.Item($row, 4)
        if ($dIsNumericLooking) {
            # Force text storage so the digit-dot string is not coerced to a Number,
            # matching the original inlineStr (text) cell type; then restore the
            # unstyled "Normal" style so no stray s="n" attribute is left behind.
            $dCell.NumberFormat = "@"
            $dCell.Value = $dVal
            $dCell.Style = "Normal"
        } else {
            $dCell.Value = $dVal
        }
    }
    if ($eVal -ne $null) {
        $ws.Cells.Item($row, 5).Value = $eVal
    }
}
